$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1 / C1 headers: bump "10.b.1" -> "10.b.1.1" in the Kyrgyz and English titles.
$ws.Range("A1").Value = "10.b.1.1 Агымдардын түрлөрү жана алуучу өлкөлөр жана донор-өлкөлөр боюнча бөлунүшүндөгү  өнүктүрүү максатында ресурстар агымынын жалпы көлөмү"
$ws.Range("C1").Value = "10.b.1.1 Total resource flows for development, by recipient and donor countries and type of flow (e.g. official development assistance, foreign direct investment and other flows)"

# Restore the selected cell reported in the saved view state.
$ws.Range("L8").Select()
